$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 754:755, pushing the existing rows 754-804
# down to 756-806 (preserving all their data/formatting untouched).
$ws.Range("754:755").Insert()

# Populate the two newly inserted rows with the new daily price entries.
# Row 754
$ws.Range("A754").Value = 4
$ws.Range("B754").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C754").Value = "Los Lagos"
$ws.Range("D754").Value = 45021
$ws.Range("E754").Value = 10
$ws.Range("F754").Value = "Fruta"
$ws.Range("G754").Value = 100106
$ws.Range("H754").Value = "Oleaginosos"
$ws.Range("I754").Value = 100106002
$ws.Range("J754").Value = "Palta"
$ws.Range("K754").Value = "Hass"
$ws.Range("L754").Value = "Primera"
$ws.Range("M754").Value = 100
$ws.Range("N754").Value = 5700
$ws.Range("O754").Value = 5700
$ws.Range("P754").Value = 5700
$ws.Range("Q754").Value = "$/kilo (en caja de 17 kilos)"
$ws.Range("R754").Value = "Provincia de Quillota"
$ws.Range("S754").Value = 5700
$ws.Range("T754").Value = 1

# Row 755
$ws.Range("A755").Value = 4
$ws.Range("B755").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C755").Value = "Los Lagos"
$ws.Range("D755").Value = 45021
$ws.Range("E755").Value = 10
$ws.Range("F755").Value = "Fruta"
$ws.Range("G755").Value = 100106
$ws.Range("H755").Value = "Oleaginosos"
$ws.Range("I755").Value = 100106002
$ws.Range("J755").Value = "Palta"
$ws.Range("K755").Value = "Hass"
$ws.Range("L755").Value = "Segunda"
$ws.Range("M755").Value = 100
$ws.Range("N755").Value = 5200
$ws.Range("O755").Value = 5200
$ws.Range("P755").Value = 5200
$ws.Range("Q755").Value = "$/kilo (en caja de 17 kilos)"
$ws.Range("R755").Value = "Provincia de Quillota"
$ws.Range("S755").Value = 5200
$ws.Range("T755").Value = 1
